# Helper: find a shape on a slide by its persisted shape Id (robust to
# shape-name collisions, e.g. several "Google Shape;245;p14" shapes per slide).
function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$p = $ppt.ActivePresentation

# --- Slide 3 : drop the trailing period from the legislation blurb --------
$s3 = $p.Slides.Item(3)
$sh3 = Get-ShapeById $s3 225
$sh3.TextFrame.TextRange.Text = "Legislation recently passed regarding hearing aids that will impact a large percent of the population"

# --- Slide 4 : collapse the 3 runs "...in " + "jupyter" + " notebook" -----
#     into a single run "...in Jupyter Notebook"
$s4 = $p.Slides.Item(4)
$sh4 = Get-ShapeById $s4 3
$tr4 = $sh4.TextFrame.TextRange
$para1 = $tr4.Paragraphs(1, 1)
$para1.Text = "Using Pandas, the Python library, in Jupyter Notebook"

# --- Slide 5 : extend the S3 Bucket description ---------------------------
$s5 = $p.Slides.Item(5)
$sh5 = Get-ShapeById $s5 244
$tr5 = $sh5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(1, 1)
$para5.Text = "S3 Bucket, an Amazon Web Service application for data storage"

# --- Slide 6 : split the ML blurb into two paragraphs ---------------------
$s6 = $p.Slides.Item(6)
$sh6 = Get-ShapeById $s6 3
$tr6 = $sh6.TextFrame.TextRange
$tr6.Text = "Using competitive model structure and training on Costco locations to recommend where Costco should open a new location, specifically for hearing centers. `rUsing precision and recall rather than accuracy to predict where they should have a hearing center but don't."

# The former single, centre-forced paragraph becomes two left (default)
# aligned paragraphs once the text is split.
$paras6 = $tr6.Paragraphs()
for ($i = 1; $i -le $paras6.Count; $i++) {
    $paras6.Item($i).ParagraphFormat.Alignment = 1
}
